$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.318.43'
$ws.Range('E2').Value = '  -0.97%  '

$ws.Range('D3').Value = '3.575.46'
$ws.Range('E3').Value = '  +2.35%  '

$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').Value = '607.89'
$ws.Range('E5').Value = '  +0.25%  '

$ws.Range('D6').Value = '144.88'
$ws.Range('E6').Value = '  -0.45%  '

$ws.Range('D7').Value = '3.574.40'
$ws.Range('E7').Value = '  +2.35%  '

$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('E9').Value = '  +1.53%  '

$ws.Range('E10').Value = '  -2.82%  '

$ws.Range('D11').Value = '7.98'
$ws.Range('E11').Value = '  +0.28%  '

$ws.Range('E12').Value = '  -1.31%  '

$ws.Range('D13').Value = '4.182.60'
$ws.Range('E13').Value = '  +2.40%  '

$ws.Range('E14').Value = '  -1.91%  '

$ws.Range('D15').Value = '30.28'
$ws.Range('E15').Value = '  -2.77%  '

$ws.Range('D16').Value = '3.574.22'
$ws.Range('E16').Value = '  +2.27%  '

$ws.Range('D17').Value = '66.394.89'
$ws.Range('E17').Value = '  -1.05%  '

$ws.Range('D18').Value = '11.71'
$ws.Range('E18').Value = '  +8.91%  '

$ws.Range('E19').Value = '  -1.36%  '

$ws.Range('D20').Value = '6.22'
$ws.Range('E20').Value = '  -1.09%  '

$ws.Range('D21').Value = '14.98'
$ws.Range('E21').Value = '  -2.45%  '

$ws.Range('D22').Value = '431.40'
$ws.Range('E22').Value = '  +0.62%  '

$ws.Range('E23').Value = '  +1.06%  '

$ws.Range('D24').Value = '78.66'
$ws.Range('E24').Value = '  -0.52%  '

$ws.Range('D25').Value = '3.717.87'
$ws.Range('E25').Value = '  +2.51%  '

$ws.Range('E26').Value = '  -0.05%  '

$ws.Range('D27').Value = '0.0000121'
$ws.Range('E27').Value = '  +3.47%  '

$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').Value = '2.53'
$ws.Range('E28').Value = '  +0.88%  '

$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '8.06'
$ws.Range('E29').Value = '  -0.96%  '

$ws.Range('D30').Value = '9.26'
$ws.Range('E30').Value = '  -4.68%  '

$ws.Range('E31').Value = '  +0.04%  '

$ws.Range('D32').Value = '1.49'
$ws.Range('E32').Value = '  -3.97%  '

$ws.Range('D33').Value = '0.159'
$ws.Range('E33').Value = '  -3.82%  '

$ws.Range('D34').Value = '3.573.72'
$ws.Range('E34').Value = '  +2.46%  '

$ws.Range('D35').Value = '25.48'
$ws.Range('E35').Value = '  +0.55%  '

$ws.Range('B36').Value = 'USDe'
$ws.Range('C36').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.04%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '1.76'
$ws.Range('E37').Value = '  -0.78%  '

$ws.Range('D38').Value = '7.91'
$ws.Range('E38').Value = '  -0.32%  '

$ws.Range('D39').Value = '5.66'
$ws.Range('E39').Value = '  -1.12%  '

$ws.Range('E40').Value = '  +0.08%  '

$ws.Range('D41').Value = '172.11'
$ws.Range('E41').Value = '  -1.74%  '

$ws.Range('D42').Value = '0.0859'
$ws.Range('E42').Value = '  -3.75%  '

$ws.Range('D43').Value = '5.30'
$ws.Range('E43').Value = '  -0.50%  '

$ws.Range('E44').Value = '  +0.78%  '

$ws.Range('E45').Value = '  -2.78%  '

$ws.Range('D46').Value = '45.83'
$ws.Range('E46').Value = '  -0.86%  '

$ws.Range('D47').Value = '1.22'
$ws.Range('E47').Value = '  +1.02%  '

$ws.Range('D48').Value = '26.02'
$ws.Range('E48').Value = '  -6.87%  '

$ws.Range('D49').Value = '2.40'
$ws.Range('E49').Value = '  +0.06%  '

$ws.Range('D50').Value = '7.15'
$ws.Range('E50').Value = '  -2.45%  '

$ws.Range('D51').Value = '0.952'
$ws.Range('E51').Value = '  -2.40%  '
